$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, "Nguyễn Xuân Đạt", "2025-12-09", "08:00:00", "17:00:00", "Có mặt"),
    @(5, "Trần Thị Ánh Nhi", "2025-12-09", "08:00:00", "17:00:00", "Có mặt"),
    @(6, "Nguyễn Quang Hoài Đức", "2025-12-09", "08:00:00", "17:00:00", "Có mặt"),
    @(7, "Bùi Anh Dũng", "2025-12-09", "08:00:00", "17:00:00", "Có mặt"),
    @(8, "Bùi Anh Dũng", "2025-12-10", "08:00:00", "17:00:00", "Có mặt"),
    @(9, "Bùi Anh Dũng", "2025-12-07", "08:00:00", "17:00:00", "Có mặt")
)

$rowIndex = 4
foreach ($row in $data) {
    # Force Text number format on the row's B:F cells first so that
    # date-/time-looking strings are stored as literal text, not
    # auto-converted into numeric date/time serials.
    $rangeAddr = "B" + $rowIndex + ":F" + $rowIndex
    $ws.Range($rangeAddr).NumberFormat = "@"

    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}
